# Generate Report for Handoff
# - Status moves from "In Translation" to "Ready for handoff" (Overview + per-locale sheets)
# - Latest HO Xliff Generate Date / Latest Handoff Datetime timestamps bump forward
# - Status/date columns widen to fit the new values

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Status: "In Translation" -> "Ready for handoff" ----------------------
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("C2").Value = "Ready for handoff"

# --- Timestamps -------------------------------------------------------------
# Overview's "Latest HO Xliff Generate Date" and de-de's "Latest Handoff Datetime"
# shared the same value before the edit and still match after it.
$wsOverview.Range("G2").Value = "2016-08-26 09:04:13"
$wsDeDe.Range("H2").Value = "2016-08-26 09:04:13"

# zh-cn's "Latest Handoff Datetime" bumps independently.
$wsZhCn.Range("H2").Value = "2016-08-26 09:04:02"

# --- Column widths ------------------------------------------------------
# Overview columns E/F (zh-cn / de-de status) and the per-locale "Status" column
# widen to accommodate "Ready for handoff".
$wsOverview.Columns.Item(5).ColumnWidth = 17.2159881591797
$wsOverview.Columns.Item(6).ColumnWidth = 17.2159881591797
$wsZhCn.Columns.Item(3).ColumnWidth = 17.2159881591797
$wsDeDe.Columns.Item(3).ColumnWidth = 17.2159881591797
